$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 45982
$ws.Range("D8").Value = 172.56
$ws.Range("E8").Value = 161.93
$ws.Range("F8").Value = 171.93
$ws.Range("G8").Value = 162.09
$ws.Range("A9").Value = 45982
$ws.Range("D9").Value = 172.56
$ws.Range("E9").Value = 161.93
$ws.Range("F9").Value = 171.93
$ws.Range("G9").Value = 162.09
$ws.Range("A10").Value = 45982
$ws.Range("D10").Value = 174.61
$ws.Range("E10").Value = 164.6
$ws.Range("F10").Value = 174.6
$ws.Range("G10").Value = 165.09
$ws.Range("A11").Value = 45981
$ws.Range("D11").Value = 172.2
$ws.Range("E11").Value = 162.07
$ws.Range("F11").Value = 172.07
$ws.Range("G11").Value = 162.22999999999999
$ws.Range("A12").Value = 45981
$ws.Range("D12").Value = 172.2
$ws.Range("E12").Value = 162.07
$ws.Range("F12").Value = 172.07
$ws.Range("G12").Value = 162.22999999999999
$ws.Range("A13").Value = 45981
$ws.Range("D13").Value = 174.24
$ws.Range("E13").Value = 165.06
$ws.Range("F13").Value = 175.06
$ws.Range("G13").Value = 165.55
$ws.Range("A17").Value = 45982
$ws.Range("D17").Value = 177.91
$ws.Range("E17").Value = 167.24
$ws.Range("F17").Value = 177.24
$ws.Range("A18").Value = 45981
$ws.Range("D18").Value = 177.55
$ws.Range("E18").Value = 167.83
$ws.Range("F18").Value = 177.83
$ws.Range("A22").Value = 45982
$ws.Range("D22").Value = 173.58
$ws.Range("E22").Value = 163.66999999999999
$ws.Range("F22").Value = 173.27
$ws.Range("G22").Value = 164.95
$ws.Range("A23").Value = 45982
$ws.Range("D23").Value = 179.39
$ws.Range("E23").Value = 168.42
$ws.Range("F23").Value = 178.42
$ws.Range("A24").Value = 45982
$ws.Range("D24").Value = 179.19
$ws.Range("E24").Value = 168.79
$ws.Range("F24").Value = 178.79
$ws.Range("A25").Value = 45982
$ws.Range("D25").Value = 180.01
$ws.Range("E25").Value = 168.2
$ws.Range("F25").Value = 178.2
$ws.Range("G25").Value = 168.24
$ws.Range("A26").Value = 45982
$ws.Range("D26").Value = 178.7
$ws.Range("E26").Value = 169.66
$ws.Range("F26").Value = 179.66
$ws.Range("A27").Value = 45981
$ws.Range("D27").Value = 173.11
$ws.Range("E27").Value = 164.01
$ws.Range("F27").Value = 173.61
$ws.Range("G27").Value = 165.3
$ws.Range("A28").Value = 45981
$ws.Range("D28").Value = 179.02
$ws.Range("E28").Value = 168.55
$ws.Range("F28").Value = 178.55
$ws.Range("A29").Value = 45981
$ws.Range("D29").Value = 178.82
$ws.Range("E29").Value = 168.8
$ws.Range("F29").Value = 178.8
$ws.Range("A30").Value = 45981
$ws.Range("D30").Value = 179.64
$ws.Range("E30").Value = 168.21
$ws.Range("F30").Value = 178.21
$ws.Range("G30").Value = 168.25
$ws.Range("A31").Value = 45981
$ws.Range("D31").Value = 178.34
$ws.Range("E31").Value = 169.78
$ws.Range("F31").Value = 179.78
$ws.Range("A35").Value = 45982
$ws.Range("D35").Value = 172.85
$ws.Range("E35").Value = 161.59
$ws.Range("F35").Value = 170.59
$ws.Range("A36").Value = 45981
$ws.Range("D36").Value = 172.49
$ws.Range("E36").Value = 162.05000000000001
$ws.Range("F36").Value = 171.05
$ws.Range("A40").Value = 45982
$ws.Range("D40").Value = 178.62
$ws.Range("E40").Value = 167.2
$ws.Range("F40").Value = 177.2
$ws.Range("A41").Value = 45982
$ws.Range("D41").Value = 178.32
$ws.Range("E41").Value = 167.62
$ws.Range("F41").Value = 177.62
$ws.Range("A42").Value = 45981
$ws.Range("D42").Value = 178.24
$ws.Range("E42").Value = 167.76
$ws.Range("F42").Value = 177.76
$ws.Range("A43").Value = 45981
$ws.Range("D43").Value = 177.94
$ws.Range("E43").Value = 168.18
$ws.Range("F43").Value = 178.18
$ws.Range("A47").Value = 45982
$ws.Range("D47").Value = 172.33
$ws.Range("E47").Value = 163.27000000000001
$ws.Range("F47").Value = 173.27
$ws.Range("A48").Value = 45982
$ws.Range("D48").Value = 172.29
$ws.Range("E48").Value = 163.43
$ws.Range("F48").Value = 173.43
$ws.Range("A49").Value = 45981
$ws.Range("D49").Value = 171.88
$ws.Range("E49").Value = 163.43
$ws.Range("F49").Value = 173.43
$ws.Range("A50").Value = 45981
$ws.Range("D50").Value = 171.84
$ws.Range("E50").Value = 163.59
$ws.Range("F50").Value = 173.59
$ws.Range("A54").Value = 45982
$ws.Range("D54").Value = 188.66
$ws.Range("E54").Value = 177.89
$ws.Range("F54").Value = 187.89
$ws.Range("A55").Value = 45982
$ws.Range("D55").Value = 176.36
$ws.Range("E55").Value = 174.62
$ws.Range("F55").Value = 184.62
$ws.Range("A56").Value = 45982
$ws.Range("D56").Value = 178.76
$ws.Range("A57").Value = 45982
$ws.Range("D57").Value = 178.23
$ws.Range("E57").Value = 168.88
$ws.Range("A58").Value = 45982
$ws.Range("D58").Value = 174.14
$ws.Range("E58").Value = 164.94
$ws.Range("F58").Value = 174.94
$ws.Range("A59").Value = 45982
$ws.Range("D59").Value = 180.87
$ws.Range("E59").Value = 175.84
$ws.Range("A60").Value = 45981
$ws.Range("D60").Value = 188.29
$ws.Range("E60").Value = 178.3
$ws.Range("F60").Value = 188.3
$ws.Range("A61").Value = 45981
$ws.Range("D61").Value = 175.99
$ws.Range("E61").Value = 175.08
$ws.Range("F61").Value = 185.08
$ws.Range("A62").Value = 45981
$ws.Range("D62").Value = 178.39
$ws.Range("A63").Value = 45981
$ws.Range("D63").Value = 177.87
$ws.Range("E63").Value = 169.35
$ws.Range("A64").Value = 45981
$ws.Range("D64").Value = 173.78
$ws.Range("E64").Value = 165.4
$ws.Range("F64").Value = 175.4
$ws.Range("A65").Value = 45981
$ws.Range("D65").Value = 180.52
$ws.Range("E65").Value = 176.27
